# Add 2022-Q1 sheet (feat: add 2022-Q1 data)
$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Item("2021-Q4")

# 1. Create the new "2022-Q1" sheet right before "总计" so tab order becomes
#    2021-Q2, 2021-Q3, 2021-Q4, 2022-Q1, 总计.
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# NOTE: inserting a sheet "before $total" re-seats the (positional) $total
# handle onto the newly inserted sheet itself, since sheet references here
# track a slot index rather than a stable identity. Re-resolve "总计" by
# name so later edits land on the real totals sheet, not the new one.
$total = $wb.Worksheets.Item("总计")

# Seed formatting (header style / index-column style / column layout) by
# copying the format of an existing quarter sheet that already has the
# exact same look.
$q4.Range("A1:H9").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holdings data for 2022-Q1.
$rows = @(
    @("001915", "宝盈医疗健康沪港深股票", "5.55", "93.11", "6.95", "0.3857", 3),
    @("090020", "大成健康产业混合",       "3.76", "91.73", "5.78", "0.2173", 6),
    @("012045", "大成医药健康股票A",      "2.87", "93.58", "7.26", "0.2084", 4),
    @("519673", "银河康乐股票",           "2.31", "92.35", "3.42", "0.0790", 10),
    @("012046", "大成医药健康股票C",      "0.25", "93.58", "7.26", "0.0182", 4),
    @("008884", "博远博锐混合A",          "0.19", "86.59", "6.34", "0.0120", 1),
    @("001563", "华富健康文娱灵活配置混合", "0.10", "90.86", "3.28", "0.0033", 4),
    @("008885", "博远博锐混合C",          "0.02", "86.59", "6.34", "0.0013", 1)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = ($r - 2)

    $bCell = $newSheet.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[0]

    $newSheet.Cells.Item($r, 3).Value = $row[1]

    $dCell = $newSheet.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[2]

    $eCell = $newSheet.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[3]

    $fCell = $newSheet.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[4]

    $gCell = $newSheet.Cells.Item($r, 7)
    $gCell.NumberFormat = "@"
    $gCell.Value = $row[5]

    $newSheet.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# 2. Insert the 2022-Q1 summary row at the top of the "总计" (totals) sheet,
#    pushing the existing 2021-Q4 / 2021-Q3 / 2021-Q2 rows down by one.
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 0.93

# Restore the index-column style (thin border / bold / centered) on A2 to
# match the other rows in this column.
$q4.Range("A2").Copy()
$total.Range("A2").PasteSpecial(-4122)

# Renumber the running index in column A for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
